# Append four new CCP rows (KELER / Q4-2017: TEA, KGA, TP KGA, CEEGEX KGA)
# to the end of the "all_in_one" sheet, growing the used range from
# A1:AC340 to A1:AC344.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 341
$ws.Range("A341").Value = 339
$ws.Range("B341").Value = "KELER"
$ws.Range("C341").Value = "Q4-2017"
$ws.Range("D341").Value = "TEA"
$ws.Range("E341").Value = 504653.1888824402
$ws.Range("F341").Value = 16121.75146707938
$ws.Range("G341").Value = 2579480.234732701
$ws.Range("H341").Value = 5046108.209195848
$ws.Range("I341").Value = 5046108.209195848
$ws.Range("J341").Value = 0
$ws.Range("K341").Value = 0
$ws.Range("L341").Value = 0
$ws.Range("M341").Value = 5046108.209195848
$ws.Range("N341").Value = 0
$ws.Range("O341").Value = 2755897.59554814
$ws.Range("P341").Value = 764022.3593160097
$ws.Range("Q341").Value = 4814684
$ws.Range("R341").Value = 6351926
$ws.Range("S341").Value = 63136358.12420253
$ws.Range("T341").Value = 60364677.56497066
$ws.Range("U341").Value = 11
$ws.Range("V341").Value = 12
$ws.Range("W341").Value = 0
$ws.Range("Y341").Value = 22
$ws.Range("Z341").Value = 0
$ws.Range("AA341").Value = 11932872.98360656
$ws.Range("AB341").Value = 0
$ws.Range("AC341").Value = 42351496.90609522

# Row 342
$ws.Range("A342").Value = 340
$ws.Range("B342").Value = "KELER"
$ws.Range("C342").Value = "Q4-2017"
$ws.Range("D342").Value = "KGA"
$ws.Range("E342").Value = 942879.1900432063
$ws.Range("F342").Value = 16121.75146707938
$ws.Range("G342").Value = 2579480.234732701
$ws.Range("H342").Value = 9428000.257948024
$ws.Range("I342").Value = 9428000.257948024
$ws.Range("J342").Value = 0
$ws.Range("K342").Value = 0
$ws.Range("L342").Value = 0
$ws.Range("M342").Value = 9428000.257948024
$ws.Range("N342").Value = 0
$ws.Range("O342").Value = 2455221.719561488
$ws.Range("P342").Value = 1426851.958320716
$ws.Range("Q342").Value = 17043517
$ws.Range("R342").Value = 28397506
$ws.Range("S342").Value = 65411808.06437584
$ws.Range("T342").Value = 62619609.24421229
$ws.Range("U342").Value = 8
$ws.Range("V342").Value = 7
$ws.Range("W342").Value = 0
$ws.Range("X342").Value = 875
$ws.Range("Y342").Value = 12
$ws.Range("Z342").Value = 0
$ws.Range("AA342").Value = 23738.90163934425
$ws.Range("AB342").Value = 0
$ws.Range("AC342").Value = 26811630.04297515

# Row 343
$ws.Range("A343").Value = 341
$ws.Range("B343").Value = "KELER"
$ws.Range("C343").Value = "Q4-2017"
$ws.Range("D343").Value = "TP KGA"
$ws.Range("E343").Value = 322964.6449990327
$ws.Range("F343").Value = 16121.75146707938
$ws.Range("G343").Value = 2579480.234732701
$ws.Range("H343").Value = 3229375.285355001
$ws.Range("I343").Value = 3229375.285355001
$ws.Range("J343").Value = 0
$ws.Range("K343").Value = 0
$ws.Range("L343").Value = 0
$ws.Range("M343").Value = 3229375.285355001
$ws.Range("N343").Value = 0
$ws.Range("Q343").Value = 0
$ws.Range("R343").Value = 9787349
$ws.Range("S343").Value = 20024435.17669093
$ws.Range("T343").Value = 19989468.28851486
$ws.Range("U343").Value = 0
$ws.Range("V343").Value = 38
$ws.Range("W343").Value = 0
$ws.Range("Y343").Value = 0
$ws.Range("Z343").Value = 0
$ws.Range("AA343").Value = 0
$ws.Range("AB343").Value = 0
$ws.Range("AC343").Value = 0

# Row 344
$ws.Range("A344").Value = 342
$ws.Range("B344").Value = "KELER"
$ws.Range("C344").Value = "Q4-2017"
$ws.Range("D344").Value = "CEEGEX KGA"
$ws.Range("E344").Value = 112241.4167795189
$ws.Range("F344").Value = 16121.75146707938
$ws.Range("G344").Value = 2579480.234732701
$ws.Range("H344").Value = 1111000
$ws.Range("I344").Value = 1111000
$ws.Range("J344").Value = 0
$ws.Range("K344").Value = 0
$ws.Range("L344").Value = 0
$ws.Range("M344").Value = 1111000
$ws.Range("N344").Value = 0
$ws.Range("O344").Value = 126998.4047359257
$ws.Range("P344").Value = 39534.6277334805
$ws.Range("Q344").Value = 0
$ws.Range("R344").Value = 136000
$ws.Range("S344").Value = 13961141.68599544
$ws.Range("T344").Value = 13348961.19816857
$ws.Range("U344").Value = 0
$ws.Range("V344").Value = 24
$ws.Range("W344").Value = 0
$ws.Range("Y344").Value = 0
$ws.Range("Z344").Value = 0
$ws.Range("AA344").Value = 0
$ws.Range("AB344").Value = 0
$ws.Range("AC344").Value = 0

# Reuse the existing bold/centered/bordered style from column A (same as A2:A340) for the new index cells
$ws.Range("A340").Copy()
$ws.Range("A341:A344").PasteSpecial(-4122)
$excel.CutCopyMode = $false
